$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped entirely from the data set:
# original row 26 "RM 232" and (after that shift) the row that held "SC 92".
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Toggle individual F-column (and one E-column) values that moved between
# "missing" (blank) and populated states, keyed off the final row numbers.
$ws.Range("F6").Value = 16.43
$ws.Range("F8").ClearContents()
$ws.Range("F12").Value = 17.45
$ws.Range("F14").ClearContents()
$ws.Range("F17").Value = 17.78
$ws.Range("F18").Value = 18.35
$ws.Range("F19").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("F23").Value = 16.48

$ws.Range("E27").Value = -10
$ws.Range("F27").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("E32").ClearContents()

Write-Output "edits applied"
